# Apply new championship results ("ft: new champs && invoke params")
$wb = $excel.ActiveWorkbook

$wsEquipos = $wb.Worksheets.Item("Equipos")
$wsResultados = $wb.Worksheets.Item("Resultados")

# --- Sheet "Equipos": update points (D) and localias faltantes (E) ---
$wsEquipos.Range("D2").Value = 3
$wsEquipos.Range("E2").Value = 1

$wsEquipos.Range("D3").Value = 6
$wsEquipos.Range("E3").Value = 1

$wsEquipos.Range("D4").Value = 3
$wsEquipos.Range("E4").Value = 2

$wsEquipos.Range("D5").Value = 6
$wsEquipos.Range("E5").Value = 2

# --- Sheet "Resultados": update matchups (C, D) and scores (E) ---
$wsResultados.Range("C3").Value = "D"
$wsResultados.Range("E3").Value = "4:4"

$wsResultados.Range("C4").Value = "B"
$wsResultados.Range("D4").Value = "C"
$wsResultados.Range("E4").Value = "2:4"

$wsResultados.Range("C6").Value = "A"
$wsResultados.Range("D6").Value = "C"
$wsResultados.Range("E6").Value = "3:2"

$wsResultados.Range("C7").Value = "B"
$wsResultados.Range("D7").Value = "D"
$wsResultados.Range("E7").Value = "2:7"

$wsResultados.Range("D9").Value = "B"
$wsResultados.Range("E9").Value = "4:4"

$wsResultados.Range("D10").Value = "D"
$wsResultados.Range("E10").Value = "2:3"

$wsResultados.Range("E12").Value = "5:1"

$wsResultados.Range("C13").Value = "B"
$wsResultados.Range("D13").Value = "C"
$wsResultados.Range("E13").Value = "3:1"

$wsResultados.Range("C15").Value = "A"
$wsResultados.Range("D15").Value = "C"
$wsResultados.Range("E15").Value = "2:0"

$wsResultados.Range("C16").Value = "B"
$wsResultados.Range("D16").Value = "D"
$wsResultados.Range("E16").Value = "3:4"

$wsResultados.Range("E18").Value = "0:2"

$wsResultados.Range("E19").Value = "4:2"
